$wb = $excel.ActiveWorkbook

# Sheet "OFF": update Week 13 row (row 2) target depth data
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B2").Value = 198
$wsOff.Range("C2").Value = 150
$wsOff.Range("D2").Value = 68
$wsOff.Range("E2").Value = 35

# Sheet "DEF": update Week 13 row (row 2) target depth data
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B2").Value = 212
$wsDef.Range("C2").Value = 162
$wsDef.Range("D2").Value = 60
$wsDef.Range("E2").Value = 22
$wsDef.Range("G2").Value = 3
